$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Append the new "Knärot" section (heading, body paragraphs, refs)
#    right after the last paragraph of the body
#    ("BILAGA 1 - Fridlysta arter"), i.e. still before sectPr.
#    Pass 1: insert each paragraph as PLAIN text (no italics yet) so
#    that newly-created empty paragraphs never pick up stray italic
#    formatting from a preceding run.
# ---------------------------------------------------------------------
$anchor = $d.Paragraphs.Last
$newParas = @()

# --- New paragraph 1: style=Heading1 ---
$anchor.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Style = 'Heading1'
$p1.Range.Text = 'Knärot – ekologi samt krav på livsmiljön'
$anchor = $p1
$newParas += $p1

# --- New paragraph 2: style=Normal ---
$anchor.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Style = 'Normal'
$p2.Range.Text = 'Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).'
$anchor = $p2
$newParas += $p2

# --- New paragraph 3: style=Normal ---
$anchor.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Style = 'Normal'
$p3.Range.Text = 'Samuel Johnsons doktorsavhandling “Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“ (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: “Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” Vidare “More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”'
$anchor = $p3
$newParas += $p3

# --- New paragraph 4: style=Normal ---
$anchor.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$p4.Style = 'Normal'
$p4.Range.Text = 'Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: “In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”'
$anchor = $p4
$newParas += $p4

# --- New paragraph 5: style=Normal ---
$anchor.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Last
$p5.Style = 'Normal'
$p5.Range.Text = 'En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).'
$anchor = $p5
$newParas += $p5

# --- New paragraph 6: style=Normal ---
$anchor.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Last
$p6.Style = 'Normal'
$p6.Range.Text = 'Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).'
$anchor = $p6
$newParas += $p6

# --- New paragraph 7: style=Heading2 ---
$anchor.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs.Last
$p7.Style = 'Heading2'
$p7.Range.Text = 'Referenser - knärot'
$anchor = $p7
$newParas += $p7

# --- New paragraph 8: style=Normal ---
$anchor.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs.Last
$p8.Style = 'Normal'
$p8.Range.Text = 'de Graaf M & Roberts M.R., 2009. Short-term response of the herbaceous layer within leave patches after harvest. Forest Ecology and Management 257, 1014-1025'
$anchor = $p8
$newParas += $p8

# --- New paragraph 9: style=Normal ---
$anchor.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs.Last
$p9.Style = 'Normal'
$p9.Range.Text = 'Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. Ecological Applications, 22, 2049-2064 '
$anchor = $p9
$newParas += $p9

# --- New paragraph 10: style=Normal ---
$anchor.Range.InsertParagraphAfter()
$p10 = $d.Paragraphs.Last
$p10.Style = 'Normal'
$p10.Range.Text = 'Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. Interactive effects of drought and edge exposure on old-growth forest understory species. Landscape Ecology, 37, sid 1839-1853'
$anchor = $p10
$newParas += $p10

# --- New paragraph 11: style=Normal ---
$anchor.Range.InsertParagraphAfter()
$p11 = $d.Paragraphs.Last
$p11.Style = 'Normal'
$p11.Range.Text = 'Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. Biological legacies buffer local species extinction after logging. Journal of Applied Ecology. 51, 53-62.'
$anchor = $p11
$newParas += $p11

# --- New paragraph 12: style=Normal ---
$anchor.Range.InsertParagraphAfter()
$p12 = $d.Paragraphs.Last
$p12.Style = 'Normal'
$p12.Range.Text = 'Skogsstyrelsen, 2022. Vägledning för hänsyn till knärot. https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/'
$anchor = $p12
$newParas += $p12

# --- New paragraph 13: style=Normal ---
$anchor.Range.InsertParagraphAfter()
$p13 = $d.Paragraphs.Last
$p13.Style = 'Normal'
$p13.Range.Text = 'SLU Artdatabanken, 2021. Artfaktablad. Naturvård – artfakta. SLU Artdatabanken, Uppsala '
$anchor = $p13
$newParas += $p13

# ---------------------------------------------------------------------
# 2. Pass 2: re-find each italic segment *within its own paragraph*
#    and flip Font.Italic on, left to right, so repeated text fragments
#    resolve to the correct occurrence.
# ---------------------------------------------------------------------
# Paragraph 3 italics
$searchStart = $p3.Range.Start
$paraEnd = $p3.Range.End
$sr2 = $d.Range($searchStart, $paraEnd)
$sr2.Find.Execute('Samuel Johnsons doktorsavhandling ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchStart = $sr2.End
$sr = $d.Range($searchStart, $paraEnd)
$sr.Find.Execute('“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sr.Font.Italic = $true
$searchStart = $sr.End
$sr2 = $d.Range($searchStart, $paraEnd)
$sr2.Find.Execute(' (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchStart = $sr2.End
$sr = $d.Range($searchStart, $paraEnd)
$sr.Find.Execute('“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sr.Font.Italic = $true
$searchStart = $sr.End
$sr2 = $d.Range($searchStart, $paraEnd)
$sr2.Find.Execute('Vidare ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchStart = $sr2.End
$sr = $d.Range($searchStart, $paraEnd)
$sr.Find.Execute('“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sr.Font.Italic = $true
$searchStart = $sr.End

# Paragraph 4 italics
$searchStart = $p4.Range.Start
$paraEnd = $p4.Range.End
$sr2 = $d.Range($searchStart, $paraEnd)
$sr2.Find.Execute('Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchStart = $sr2.End
$sr = $d.Range($searchStart, $paraEnd)
$sr.Find.Execute('“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sr.Font.Italic = $true
$searchStart = $sr.End

# Paragraph 8 italics
$searchStart = $p8.Range.Start
$paraEnd = $p8.Range.End
$sr2 = $d.Range($searchStart, $paraEnd)
$sr2.Find.Execute('de Graaf M & Roberts M.R., 2009. ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchStart = $sr2.End
$sr = $d.Range($searchStart, $paraEnd)
$sr.Find.Execute('Short-term response of the herbaceous layer within leave patches after harvest. ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sr.Font.Italic = $true
$searchStart = $sr.End
$sr2 = $d.Range($searchStart, $paraEnd)
$sr2.Find.Execute('Forest Ecology and Management 257, 1014-1025', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchStart = $sr2.End

# Paragraph 9 italics
$searchStart = $p9.Range.Start
$paraEnd = $p9.Range.End
$sr2 = $d.Range($searchStart, $paraEnd)
$sr2.Find.Execute('Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchStart = $sr2.End
$sr = $d.Range($searchStart, $paraEnd)
$sr.Find.Execute('Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sr.Font.Italic = $true
$searchStart = $sr.End
$sr2 = $d.Range($searchStart, $paraEnd)
$sr2.Find.Execute('Ecological Applications, 22, 2049-2064 ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchStart = $sr2.End

# Paragraph 10 italics
$searchStart = $p10.Range.Start
$paraEnd = $p10.Range.End
$sr2 = $d.Range($searchStart, $paraEnd)
$sr2.Find.Execute('Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchStart = $sr2.End
$sr = $d.Range($searchStart, $paraEnd)
$sr.Find.Execute('Interactive effects of drought and edge exposure on old-growth forest understory species. ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sr.Font.Italic = $true
$searchStart = $sr.End
$sr2 = $d.Range($searchStart, $paraEnd)
$sr2.Find.Execute('Landscape Ecology, 37, sid 1839-1853', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchStart = $sr2.End

# Paragraph 11 italics
$searchStart = $p11.Range.Start
$paraEnd = $p11.Range.End
$sr2 = $d.Range($searchStart, $paraEnd)
$sr2.Find.Execute('Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchStart = $sr2.End
$sr = $d.Range($searchStart, $paraEnd)
$sr.Find.Execute('Biological legacies buffer local species extinction after logging. ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sr.Font.Italic = $true
$searchStart = $sr.End
$sr2 = $d.Range($searchStart, $paraEnd)
$sr2.Find.Execute('Journal of Applied Ecology. 51, 53-62.', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchStart = $sr2.End

# Paragraph 12 italics
$searchStart = $p12.Range.Start
$paraEnd = $p12.Range.End
$sr2 = $d.Range($searchStart, $paraEnd)
$sr2.Find.Execute('Skogsstyrelsen, 2022. ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchStart = $sr2.End
$sr = $d.Range($searchStart, $paraEnd)
$sr.Find.Execute('Vägledning för hänsyn till knärot. ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sr.Font.Italic = $true
$searchStart = $sr.End
$sr2 = $d.Range($searchStart, $paraEnd)
$sr2.Find.Execute('https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchStart = $sr2.End

# Paragraph 13 italics
$searchStart = $p13.Range.Start
$paraEnd = $p13.Range.End
$sr2 = $d.Range($searchStart, $paraEnd)
$sr2.Find.Execute('SLU Artdatabanken, 2021. ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchStart = $sr2.End
$sr = $d.Range($searchStart, $paraEnd)
$sr.Find.Execute('Artfaktablad. Naturvård – artfakta. ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sr.Font.Italic = $true
$searchStart = $sr.End
$sr2 = $d.Range($searchStart, $paraEnd)
$sr2.Find.Execute('SLU Artdatabanken, Uppsala ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchStart = $sr2.End

# ---------------------------------------------------------------------
# 3. Update the date stamp in the header: 2023-09-13 -> 2023-09-15
#    (walk every story range so headers/footers are reached too)
# ---------------------------------------------------------------------
foreach ($story in $d.StoryRanges) {
    $r = $story
    while ($r -ne $null) {
        $r.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null
        $r = $r.NextStoryRange
    }
}

Write-Host "Done"